$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos list values (prices & volume deltas) per upstream diff.
# Row 40-42 coin order also changed (ARBITRUM, Stacks, TheGraph reshuffled).

$ws.Range("D2").Value = "51.501.95"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").Value = "3.021.29"
$ws.Range("E3").Value = "  +3.46%  "
$ws.Range("E4").Value = "  +0.10%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "379.74"
$r.Style = "Normal"
$ws.Range("E5").Value = "  +1.25%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "102.74"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +3.10%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.544"
$r.Style = "Normal"
$ws.Range("E7").Value = "  +1.82%  "
$ws.Range("E8").Value = "  +0.01%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.593"
$r.Style = "Normal"
$ws.Range("E9").Value = "  +3.76%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "36.67"
$r.Style = "Normal"
$ws.Range("E10").Value = "  +3.10%  "
$ws.Range("E11").Value = "  +0.00%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.0856"
$r.Style = "Normal"
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("D13").Value = "3.501.54"
$ws.Range("E13").Value = "  +3.65%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "18.47"
$r.Style = "Normal"
$ws.Range("E14").Value = "  +2.88%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "7.73"
$r.Style = "Normal"
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "3.025.14"
$ws.Range("E16").Value = "  +3.92%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "0.985"
$r.Style = "Normal"
$ws.Range("E17").Value = "  -0.16%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "10.59"
$r.Style = "Normal"
$ws.Range("E18").Value = "  -12.23%  "
$ws.Range("D19").Value = "51.546.59"
$ws.Range("E19").Value = "  +1.65%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "3.05"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +1.94%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "12.43"
$r.Style = "Normal"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("D22").Value = "0.0₃0961"
$ws.Range("E22").Value = "  +2.11%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "70.10"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +1.21%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "267.87"
$r.Style = "Normal"
$ws.Range("E24").Value = "  +0.92%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "3.16"
$r.Style = "Normal"
$ws.Range("E25").Value = "  +0.05%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "8.21"
$r.Style = "Normal"
$ws.Range("E26").Value = "  +5.00%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "7.46"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +6.03%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "0.170"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +5.63%  "
$ws.Range("E29").Value = "  -0.04%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "26.19"
$r.Style = "Normal"
$ws.Range("E30").Value = "  +3.49%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "0.109"
$r.Style = "Normal"
$ws.Range("E31").Value = "  +2.25%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "10.30"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +3.69%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "34.14"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +3.17%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "50.58"
$r.Style = "Normal"
$ws.Range("E34").Value = "  +0.63%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "2.06"
$r.Style = "Normal"
$ws.Range("E35").Value = "  +0.58%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "0.0450"
$r.Style = "Normal"
$ws.Range("E36").Value = "  +5.22%  "
$ws.Range("E37").Value = "  +0.01%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "3.27"
$r.Style = "Normal"
$ws.Range("E38").Value = "  +7.09%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "17.23"
$r.Style = "Normal"
$ws.Range("E39").Value = "  +5.97%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "1.86"
$r.Style = "Normal"
$ws.Range("E40").Value = "  +4.36%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "2.58"
$r.Style = "Normal"
$ws.Range("E41").Value = "  +7.39%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.280"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +9.21%  "
$ws.Range("E43").Value = "  +0.81%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "126.44"
$r.Style = "Normal"
$ws.Range("E44").Value = "  +2.83%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "3.73"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +11.19%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "22.01"
$r.Style = "Normal"
$ws.Range("E46").Value = "  +5.48%  "
$ws.Range("E47").Value = "  +0.62%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "2.38"
$r.Style = "Normal"
$ws.Range("E48").Value = "  +2.09%  "
$ws.Range("D49").Value = "2.028.48"
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("D50").Value = "3.322.91"
$ws.Range("E50").Value = "  +3.50%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.0320"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +2.50%  "
